$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F25").Value = "application instructions || env warning - species || env warning - water"
$ws.Range("F27").Value = "32_physical_and_chemical_hazards"
$ws.Range("F31").Value = "135_product_information"
$ws.Range("F32").Value = "application instructions"
$ws.Range("F33").Value = "application instructions"
$ws.Range("F34").Value = "application instructions"
$ws.Range("F35").Value = "application instructions"
$ws.Range("F37").Value = "use restrictions || application instructions"
$ws.Range("F38").Value = "use restrictions || application instructions"
$ws.Range("F39").Value = "application instructions"
$ws.Range("F41").Value = "application instructions"
$ws.Range("F42").Value = "application instructions"
$ws.Range("F44").Value = "application instructions"
$ws.Range("F47").Value = "154_pesticide_storage"
